$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15
$ws.Cells.Item($row, 1).Value = "24/10/2025"
$ws.Cells.Item($row, 2).Value = "Csikszereda M. Ciuc"
$ws.Cells.Item($row, 3).Value = 1
$ws.Cells.Item($row, 4).Value = 1
$ws.Cells.Item($row, 5).Value = "Petrolul"
$ws.Cells.Item($row, 6).Value = "D"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = 1
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 2.31
$ws.Cells.Item($row, 12).Value = 1.79
$ws.Cells.Item($row, 13).Value = 16
$ws.Cells.Item($row, 14).Value = 11
$ws.Cells.Item($row, 15).Value = 3
$ws.Cells.Item($row, 16).Value = 4
